$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column H width (matching new col definition width="10.5546875")
$ws.Columns.Item(8).ColumnWidth = 10.5546875

# New row 17: Meeting6 / All Member (same style as other task rows, bold font)
$ws.Range("A17").Value = "Meeting6"
$ws.Range("H17").Value = "All Member"

# New row 18: Set Condition with DOM and Frontend / 夏义
$ws.Range("A18").Value = "Set Condition with DOM and Frontend"
$ws.Range("H18").Value = "夏义"

# Match the bold formatting used by the other task rows (A4:A16, etc.)
$ws.Range("A17:A18").Font.Bold = $true
$ws.Range("H17:H18").Font.Bold = $true

# Update selection to match diff
$ws.Range("I18").Select() | Out-Null
